{"js": "// Fix Supreme Court mention missing from short data engineering resume.\n// 1) Neutralize summary language: \"all Black and Asian-American voters\" -> \"50M voters\"\n// 2) Replace the 4-bullet \"KEY ACHIEVEMENTS AND IMPACT / Impact\" list with the new 6-bullet list.\n\nconst BOLD_COLOR = \"#2C3E50\";\n\n// ---------------------------------------------------------------\n// Step 1: Professional summary neutral-language fix (paragraph #3)\n// ---------------------------------------------------------------\nconst body = context.document.body;\nlet paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst summaryPara = paras.items[3];\nconst summaryRange = summaryPara.getRange();\nconst summaryHits = summaryRange.search(\n  \"affecting all Black and Asian-American voters, developed\",\n  { matchCase: true }\n);\nsummaryHits.load(\"items\");\nawait context.sync();\n\nif (summaryHits.items.length > 0) {\n  summaryHits.items[0].insertText(\n    \"affecting 50M voters, developed\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// Step 2: Rebuild the achievements bullet list.\n// ---------------------------------------------------------------\nparas = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the old bullets by their exact current text so this is robust\n// even if paragraph indices differ slightly.\nconst oldTexts = {\n  discovered:\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n  algorithm: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n  redistricting: \"\u2022 Built redistricting platform used by thousands of analysts nationwide\",\n  achieved: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n};\n\nlet discoveredPara = null;\nlet algorithmPara = null;\nlet redistrictingPara = null;\nlet achievedPara = null;\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t === oldTexts.discovered) discoveredPara = paras.items[i];\n  else if (t === oldTexts.algorithm) algorithmPara = paras.items[i];\n  else if (t === oldTexts.redistricting) redistrictingPara = paras.items[i];\n  else if (t === oldTexts.achieved) achievedPara = paras.items[i];\n}\n\nif (!discoveredPara || !algorithmPara || !redistrictingPara || !achievedPara) {\n  throw new Error(\"Could not locate all four existing achievement bullets.\");\n}\n\n// Anchor all new paragraphs immediately before the first old bullet so the\n// new ones inherit plain (\"Normal\") formatting with no style bleed.\nconst anchor = discoveredPara;\n\n// Helper: insert a brand-new empty paragraph right before `anchor`, then\n// build its content via sequential \"End\" inserts so each run gets its own\n// formatting (mirrors how Word records mixed-formatting runs).\nfunction newBulletBefore(anchorPara) {\n  return anchorPara.insertParagraph(\"\", \"Before\");\n}\n\n// --- Bullet 1: Algorithmic innovation ... **73.5%** -------------------\nconst p1 = newBulletBefore(anchor);\np1.insertText(\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **\",\n  \"End\"\n);\nawait context.sync();\nconst p1Bold = p1.insertText(\"73.5%\", \"End\");\np1Bold.font.bold = true;\np1Bold.font.color = BOLD_COLOR;\nawait context.sync();\np1.insertText(\"**\", \"End\");\nawait context.sync();\n\n// --- Bullet 2: **$4.7M** savings enabled nonprofit access -------------\nconst p2 = newBulletBefore(anchor);\np2.insertText(\"\u2022 **\", \"End\");\nawait context.sync();\nconst p2Bold = p2.insertText(\"$4.7M\", \"End\");\np2Bold.font.bold = true;\np2Bold.font.color = BOLD_COLOR;\nawait context.sync();\np2.insertText(\"** savings enabled nonprofit access\", \"End\");\nawait context.sync();\n\n// --- Bullet 3: Legal precedent -----------------------------------------\nconst p3 = newBulletBefore(anchor);\np3.insertText(\"\u2022 Legal precedent: Data analysis utilized in Supreme Court case\", \"End\");\nawait context.sync();\n\n// --- Bullet 4: Expert methodology ---------------------------------------\nconst p4 = newBulletBefore(anchor);\np4.insertText(\"\u2022 Expert methodology validated at highest judicial level\", \"End\");\nawait context.sync();\n\n// --- Bullet 5: Breakthrough demographic discovery -----------------------\nconst p5 = newBulletBefore(anchor);\np5.insertText(\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"End\"\n);\nawait context.sync();\n\n// --- Bullet 6: **178%** accuracy improvement -----------------------------\nconst p6 = newBulletBefore(anchor);\np6.insertText(\"\u2022 **\", \"End\");\nawait context.sync();\nconst p6Bold = p6.insertText(\"178%\", \"End\");\np6Bold.font.bold = true;\np6Bold.font.color = BOLD_COLOR;\nawait context.sync();\np6.insertText(\"** accuracy improvement in racial classification algorithms\", \"End\");\nawait context.sync();\n\n// Finally remove the four original bullet paragraphs.\ndiscoveredPara.delete();\nalgorithmPara.delete();\nredistrictingPara.delete();\nachievedPara.delete();\nawait context.sync();\n", "ps1": "# Fix Supreme Court mention missing from short data engineering resume.\n# 1) Neutralize summary language: \"all Black and Asian-American voters\" -> \"50M voters\"\n# 2) Replace the 4-bullet \"KEY ACHIEVEMENTS AND IMPACT / Impact\" list with the new 6-bullet list.\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------\n# Helper: find a paragraph's 1-based index by its exact text (trimmed\n# of the trailing paragraph-mark carriage return).\n# -----------------------------------------------------------------\nfunction Find-ParaIndex($doc, $targetText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        $t = $t.TrimEnd([char]13)\n        if ($t -eq $targetText) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# -----------------------------------------------------------------\n# Step 1: Professional summary neutral-language fix.\n# -----------------------------------------------------------------\n$summaryIdx = Find-ParaIndex $d \"Data engineering professional with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML algorithms improving classification accuracy from 23% to 64%. Built Civic Graph data warehouse processing billions of records and platforms serving thousands of analysts nationwide.\"\n\nif ($summaryIdx -gt 0) {\n    $summaryRange = $d.Paragraphs.Item($summaryIdx).Range\n    $summaryRange.Find.Execute(\"affecting all Black and Asian-American voters, developed\", $false, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters, developed\", 2) | Out-Null\n}\n\n# -----------------------------------------------------------------\n# Step 2: Rebuild the achievements bullet list.\n# -----------------------------------------------------------------\n$BOLD_COLOR = \"#2C3E50\"\n\n$discoveredText = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters\"\n$algorithmText = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n$redistrictingText = \"\u2022 Built redistricting platform used by thousands of analysts nationwide\"\n$achievedText = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n\n$discoveredIdx = Find-ParaIndex $d $discoveredText\n$algorithmIdx = Find-ParaIndex $d $algorithmText\n$redistrictingIdx = Find-ParaIndex $d $redistrictingText\n$achievedIdx = Find-ParaIndex $d $achievedText\n\nif ($discoveredIdx -lt 0 -or $algorithmIdx -lt 0 -or $redistrictingIdx -lt 0 -or $achievedIdx -lt 0) {\n    throw \"Could not locate all four existing achievement bullets.\"\n}\n\n# Helper: insert a brand-new empty paragraph immediately before the\n# paragraph currently sitting at 1-based index $beforeIdx, returning the\n# index of the freshly-created (empty) paragraph.\nfunction Insert-BlankParaBefore($doc, $beforeIdx) {\n    $doc.Paragraphs.Item($beforeIdx).Range.InsertParagraphBefore() | Out-Null\n    return $beforeIdx\n}\n\n# Helper: append plain text to the end of paragraph $idx (i.e. just\n# before its paragraph mark) and return the Range covering the newly\n# inserted text (so callers can apply character formatting to it).\nfunction Append-RunText($doc, $idx, $text) {\n    $para = $doc.Paragraphs.Item($idx)\n    $body = $doc.Range($para.Range.Start, $para.Range.End - 1)\n    $body.Collapse(0) | Out-Null   # wdCollapseEnd\n    $startPos = $body.Start\n    $body.InsertAfter($text) | Out-Null\n    $endPos = $startPos + $text.Length\n    return $doc.Range($startPos, $endPos)\n}\n\n# All new bullets are anchored immediately before the first old bullet\n# (\"discovered\") so each new paragraph inherits plain (\"Normal\")\n# formatting with no style bleed-through.\n$anchorIdx = $discoveredIdx\n\n# --- Bullet 1: Algorithmic innovation ... **73.5%** ----------------------\n$p1 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p1 \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **\" | Out-Null\n$bold1 = Append-RunText $d $p1 \"73.5%\"\n$bold1.Font.Bold = 1\n$bold1.Font.Color = $BOLD_COLOR\nAppend-RunText $d $p1 \"**\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# --- Bullet 2: **$4.7M** savings enabled nonprofit access ----------------\n$p2 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p2 \"\u2022 **\" | Out-Null\n$bold2 = Append-RunText $d $p2 \"`$4.7M\"\n$bold2.Font.Bold = 1\n$bold2.Font.Color = $BOLD_COLOR\nAppend-RunText $d $p2 \"** savings enabled nonprofit access\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# --- Bullet 3: Legal precedent -------------------------------------------\n$p3 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p3 \"\u2022 Legal precedent: Data analysis utilized in Supreme Court case\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# --- Bullet 4: Expert methodology ----------------------------------------\n$p4 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p4 \"\u2022 Expert methodology validated at highest judicial level\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# --- Bullet 5: Breakthrough demographic discovery ------------------------\n$p5 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p5 \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# --- Bullet 6: **178%** accuracy improvement ------------------------------\n$p6 = Insert-BlankParaBefore $d $anchorIdx\nAppend-RunText $d $p6 \"\u2022 **\" | Out-Null\n$bold6 = Append-RunText $d $p6 \"178%\"\n$bold6.Font.Bold = 1\n$bold6.Font.Color = $BOLD_COLOR\nAppend-RunText $d $p6 \"** accuracy improvement in racial classification algorithms\" | Out-Null\n$anchorIdx = $anchorIdx + 1\n\n# -----------------------------------------------------------------\n# Finally, delete the four original bullet paragraphs (re-locate them\n# by exact text since indices have shifted after all the inserts).\n# -----------------------------------------------------------------\n$oldTexts = @($discoveredText, $algorithmText, $redistrictingText, $achievedText)\nforeach ($t in $oldTexts) {\n    $idx = Find-ParaIndex $d $t\n    if ($idx -gt 0) {\n        $d.Paragraphs.Item($idx).Range.Delete() | Out-Null\n    }\n}\n"}
